# Update ticket-interest ("想去人数", column F) and lowest-price
# ("最低票价", column G) figures across the four sheets, matching the
# gh-pages data refresh generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 37199
$ws1.Range("G2").Value = 68
$ws1.Range("F4").Value = 632
$ws1.Range("F7").Value = 253
$ws1.Range("F8").Value = 464
$ws1.Range("F9").Value = 834
$ws1.Range("F10").Value = 88
$ws1.Range("F11").Value = 686
$ws1.Range("F12").Value = 524
$ws1.Range("F13").Value = 34
$ws1.Range("F14").Value = 632
$ws1.Range("F15").Value = 176
$ws1.Range("F16").Value = 465
$ws1.Range("F18").Value = 1153
$ws1.Range("G19").Value = 9.9
$ws1.Range("F20").Value = 807
$ws1.Range("F21").Value = 2485
$ws1.Range("F22").Value = 986
$ws1.Range("F23").Value = 547
$ws1.Range("F24").Value = 103
$ws1.Range("F25").Value = 1152
$ws1.Range("F27").Value = 753
$ws1.Range("F29").Value = 1144

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 330
$ws2.Range("F6").Value = 55
$ws2.Range("F8").Value = 141
$ws2.Range("F9").Value = 9

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 616

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 616
$ws4.Range("F3").Value = 37199
$ws4.Range("G3").Value = 68
$ws4.Range("F5").Value = 632
$ws4.Range("F9").Value = 254
$ws4.Range("F10").Value = 464
$ws4.Range("F12").Value = 330
$ws4.Range("F14").Value = 834
$ws4.Range("F15").Value = 88
$ws4.Range("F16").Value = 686
$ws4.Range("F17").Value = 524
$ws4.Range("F18").Value = 55
$ws4.Range("F19").Value = 34
$ws4.Range("F21").Value = 141
$ws4.Range("F22").Value = 9
$ws4.Range("F24").Value = 632
$ws4.Range("F25").Value = 176
$ws4.Range("F26").Value = 465
$ws4.Range("F28").Value = 1153
$ws4.Range("G29").Value = 9.9
$ws4.Range("F30").Value = 807
$ws4.Range("F31").Value = 2485
$ws4.Range("F32").Value = 986
$ws4.Range("F33").Value = 547
$ws4.Range("F34").Value = 103
$ws4.Range("F35").Value = 1152
$ws4.Range("F38").Value = 753
$ws4.Range("F40").Value = 1144
